$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 7334.6665
$ws.Range("I29").Value = 6000
$ws.Range("J29").Value = 8002
$ws.Range("K29").Value = 18000
$ws.Range("L29").Value = 24006
$ws.Range("M29").Value = -17719
$ws.Range("N29").Value = -24568
$ws.Range("H43").Value = 7072
$ws.Range("I43").Value = 1114
$ws.Range("K43").Value = 1114
$ws.Range("M43").Value = -1045
$ws.Range("H55").Value = 432.51352
$ws.Range("J55").Value = 505.4737
$ws.Range("L55").Value = 505.4737
$ws.Range("N55").Value = -933.4737
$ws.Range("H87").Value = 89997.5
$ws.Range("J87").Value = 89997.5
$ws.Range("L87").Value = 89997.5
$ws.Range("N87").Value = -92493.5
$ws.Range("H90").Value = 89997.5
$ws.Range("J90").Value = 89997.5
$ws.Range("L90").Value = 269992.5
$ws.Range("N90").Value = -282472.5
$ws.Range("H98").Value = 1450.5883
$ws.Range("I98").Value = 728.46875
$ws.Range("J98").Value = 13004.5
$ws.Range("K98").Value = 728.46875
$ws.Range("L98").Value = 13004.5
$ws.Range("M98").Value = 769.53125
$ws.Range("N98").Value = -16000.5
$ws.Range("H116").Value = 8684
$ws.Range("I116").Value = 7141.5
$ws.Range("K116").Value = 7141.5
$ws.Range("M116").Value = -3699.5
$ws.Range("H122").Value = 1450.5883
$ws.Range("I122").Value = 728.46875
$ws.Range("J122").Value = 13004.5
$ws.Range("K122").Value = 2185.40625
$ws.Range("L122").Value = 39013.5
$ws.Range("M122").Value = 264.59375
$ws.Range("N122").Value = -43913.5
$ws.Range("H127").Value = 1181.5
$ws.Range("J127").Value = 3300
$ws.Range("L127").Value = 9900
$ws.Range("N127").Value = -19820
$ws.Range("H129").Value = 957.2143
$ws.Range("J129").Value = 2716.3333
$ws.Range("L129").Value = 8148.999899999999
$ws.Range("N129").Value = -18148.9999
$ws.Range("H132").Value = 1419.341
$ws.Range("I132").Value = 1201.2142
$ws.Range("K132").Value = 3603.6426
$ws.Range("M132").Value = -1073.6426
$ws.Range("H138").Value = 3674.7124
$ws.Range("I138").Value = 3125.2856
$ws.Range("J138").Value = 3732.9849
$ws.Range("K138").Value = 9375.856800000001
$ws.Range("L138").Value = 11198.9547
$ws.Range("M138").Value = -4235.856800000001
$ws.Range("N138").Value = -21478.9547

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4760.145
$ws.Range("I32").Value = 3453.8774
$ws.Range("J32").Value = 9683.77
$ws.Range("K32").Value = 3453.8774
$ws.Range("L32").Value = 9683.77
$ws.Range("M32").Value = -3166.8774
$ws.Range("N32").Value = -10257.77
$ws.Range("H61").Value = 7610.091
$ws.Range("I61").Value = 5633.3335
$ws.Range("J61").Value = 16505.5
$ws.Range("K61").Value = 5633.3335
$ws.Range("L61").Value = 16505.5
$ws.Range("M61").Value = -5421.3335
$ws.Range("N61").Value = -16929.5
$ws.Range("H132").Value = 2104.6606
$ws.Range("I132").Value = 1683.32
$ws.Range("J132").Value = 5615.8335
$ws.Range("K132").Value = 5049.96
$ws.Range("L132").Value = 16847.5005
$ws.Range("M132").Value = -2519.96
$ws.Range("N132").Value = -21907.5005
$ws.Range("H136").Value = 7610.091
$ws.Range("I136").Value = 5633.3335
$ws.Range("J136").Value = 16505.5
$ws.Range("K136").Value = 16900.0005
$ws.Range("L136").Value = 49516.5
$ws.Range("M136").Value = -14350.0005
$ws.Range("N136").Value = -54616.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 9489.385
$ws.Range("I26").Value = 9489.385
$ws.Range("K26").Value = 9489.385
$ws.Range("M26").Value = -9197.385
$ws.Range("H96").Value = 5713.625
$ws.Range("I96").Value = 5713.625
$ws.Range("K96").Value = 5713.625
$ws.Range("M96").Value = -2967.625
$ws.Range("H134").Value = 1968.8889
$ws.Range("I134").Value = 1468.075
$ws.Range("K134").Value = 4404.225
$ws.Range("M134").Value = -1869.225

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 318.7857
$ws.Range("J7").Value = 463.625
$ws.Range("L7").Value = 463.625
$ws.Range("N7").Value = -689.625
$ws.Range("H31").Value = 49757.78
$ws.Range("I31").Value = 2839.8
$ws.Range("K31").Value = 2839.8
$ws.Range("M31").Value = -2544.8
$ws.Range("H34").Value = 49757.78
$ws.Range("I34").Value = 2839.8
$ws.Range("K34").Value = 2839.8
$ws.Range("M34").Value = -2637.8
$ws.Range("H58").Value = 3825.3242
$ws.Range("I58").Value = 2064.0908
$ws.Range("J58").Value = 6408.467
$ws.Range("K58").Value = 2064.0908
$ws.Range("L58").Value = 6408.467
$ws.Range("M58").Value = -1861.0908
$ws.Range("N58").Value = -6814.467
$ws.Range("H103").Value = 30759.908
$ws.Range("I103").Value = 28835.9
$ws.Range("K103").Value = 28835.9
$ws.Range("M103").Value = -27663.9
$ws.Range("H131").Value = 25199
$ws.Range("I131").Value = 15296
$ws.Range("K131").Value = 15296
$ws.Range("M131").Value = -10256
$ws.Range("H134").Value = 2871.8076
$ws.Range("I134").Value = 2057.65
$ws.Range("J134").Value = 5585.6665
$ws.Range("K134").Value = 6172.950000000001
$ws.Range("L134").Value = 16756.9995
$ws.Range("M134").Value = -3637.950000000001
$ws.Range("N134").Value = -21826.9995
$ws.Range("H136").Value = 3825.3242
$ws.Range("I136").Value = 2064.0908
$ws.Range("J136").Value = 6408.467
$ws.Range("K136").Value = 6192.2724
$ws.Range("L136").Value = 19225.401
$ws.Range("M136").Value = -3642.2724
$ws.Range("N136").Value = -24325.401
$ws.Range("H141").Value = 152722.27
$ws.Range("J141").Value = 159776
$ws.Range("L141").Value = 159776
$ws.Range("N141").Value = -170136

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2223244
$ws.Range("I5").Value = 838.8182
$ws.Range("J5").Value = 8334858
$ws.Range("K5").Value = 2516.4546
$ws.Range("L5").Value = 25004574
$ws.Range("M5").Value = -2404.4546
$ws.Range("N5").Value = -25004798
$ws.Range("H33").Value = 92.09999999999999
$ws.Range("I33").Value = 94.73333
$ws.Range("K33").Value = 568.3999799999999
$ws.Range("M33").Value = -285.3999799999999
$ws.Range("H122").Value = 3035.3333
$ws.Range("J122").Value = 3502.4666
$ws.Range("L122").Value = 31522.1994
$ws.Range("N122").Value = -36422.1994
$ws.Range("H135").Value = 2223244
$ws.Range("I135").Value = 838.8182
$ws.Range("J135").Value = 8334858
$ws.Range("K135").Value = 7549.3638
$ws.Range("L135").Value = 75013722
$ws.Range("M135").Value = -5014.3638
$ws.Range("N135").Value = -75018792
$ws.Range("H136").Value = 2467.1738
$ws.Range("I136").Value = 2356.5908
$ws.Range("K136").Value = 7069.7724
$ws.Range("M136").Value = -1969.7724
$ws.Range("H137").Value = 46919.434
$ws.Range("J137").Value = 62844
$ws.Range("L137").Value = 188532
$ws.Range("N137").Value = -198732
$ws.Range("H138").Value = 1290
$ws.Range("I138").Value = 1045.7142
$ws.Range("K138").Value = 3137.1426
$ws.Range("M138").Value = 2002.8574
$ws.Range("H139").Value = 3369.8823
$ws.Range("I139").Value = 2289.5
$ws.Range("J139").Value = 5962.8
$ws.Range("K139").Value = 6868.5
$ws.Range("L139").Value = 17888.4
$ws.Range("M139").Value = -1728.5
$ws.Range("N139").Value = -28168.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3329.5103
$ws.Range("I132").Value = 2328.257
$ws.Range("K132").Value = 6984.771000000001
$ws.Range("M132").Value = -4454.771000000001
$ws.Range("H139").Value = 67969.60000000001
$ws.Range("J139").Value = 67969.60000000001
$ws.Range("L139").Value = 67969.60000000001
$ws.Range("N139").Value = -78249.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4003.6382
$ws.Range("I132").Value = 3183.8276
$ws.Range("J132").Value = 5324.4443
$ws.Range("K132").Value = 9551.4828
$ws.Range("L132").Value = 15973.3329
$ws.Range("M132").Value = -7021.4828
$ws.Range("N132").Value = -21033.3329
$ws.Range("H136").Value = 4275.485
$ws.Range("I136").Value = 2591.0833
$ws.Range("J136").Value = 8767.223
$ws.Range("K136").Value = 7773.249899999999
$ws.Range("L136").Value = 26301.669
$ws.Range("M136").Value = -5223.249899999999
$ws.Range("N136").Value = -31401.669

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2283.0605
$ws.Range("I122").Value = 1697.9642
$ws.Range("K122").Value = 5093.892599999999
$ws.Range("M122").Value = -2643.892599999999
$ws.Range("H132").Value = 1700.6383
$ws.Range("I132").Value = 1365.1708
$ws.Range("K132").Value = 4095.512400000001
$ws.Range("M132").Value = -1565.512400000001
